$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the regression output table with the new "In Class Demonstrations"
# coefficients. Values are entered column-by-column (B then C), matching the
# order the original workbook's shared strings were authored in.
$ws.Cells.Item(2, 2).Value = "0.529*"    # B2 - FFR Lag / FFR column
$ws.Cells.Item(3, 2).Value = "-3.054**"  # B3 - U Lag   / FFR column

$ws.Cells.Item(2, 3).Value = "-0.215**"  # C2 - FFR Lag / U column

# C3 ("0.697") looks like a plain number, so force it to be stored as text
# (matching the rest of the coefficient column) by temporarily switching the
# cell to a text format, then drop that number format again so the cell is
# left with its original (default) style.
$c3 = $ws.Cells.Item(3, 3)               # C3 - U Lag   / U column
$c3.NumberFormat = "@"
$c3.Value = "0.697"
$c3.Style = "Normal"

# The "Constant" and "r2_adj" rows are no longer part of the table - remove
# them entirely (shrinks the used range from A1:C5 down to A1:C3).
$ws.Range("A4:C5").Delete()
